# Implemented LinearSVM for method 1, and trained the classifier.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text: "Classifier " -> "Classifier" (trailing space removed)
$ws.Range("B3").Value = "Classifier"

# Populate row 12 with the newly-trained Linear SVM classifier metrics
$ws.Range("B12").Value = "Linear SVM"
$ws.Range("C12").Value = 27.4154318174
$ws.Range("D12").Value = 0.5965
$ws.Range("E12").Value = 0.608822103732
$ws.Range("F12").Value = 0.583376355188
$ws.Range("G12").Value = 0.590780809031
$ws.Range("H12").Value = 0.602988260406
$ws.Range("I12").Value = 0.628
$ws.Range("J12").Value = 0.565

# Match the row height Excel auto-applies once the row holds data
# (same height already used by the other populated classifier rows)
$ws.Rows.Item(12).RowHeight = 14.25

# Update the view: leave the selection on the last cell touched (J12),
# matching where the cursor ends up after entering the new row of data
$ws.Range("J12").Select()
